# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF holds a "Date" label in BF1 and a date string in BF2:BF31
# that was formatted as "2-18-2012-13" (game-log-file-name style) and
# needs to be corrected to the actual ISO-ish date "2013-02-18".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the date values are stored/kept as plain text (matching the
# source data's original inline-string type) instead of being
# auto-recognized and converted into a date serial number when written.
$dataRange = $ws.Range("BF2:BF31")
$dataRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Text -eq "2-18-2012-13") {
        $cell.Value = "2013-02-18"
    }
}
